# Converted the validity to 75%. Made all non-debug condition rely on a
# single conditions file. This script rewrites the single worksheet so
# that it only contains the two columns that are shared across all
# (non-debug) condition files: targ_right (A) and cue_valid (B), with the
# four possible combinations of 0/1 values (i.e. "75% validity" cue/target
# design collapsed to its two defining factors).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New, reduced header row.
$ws.Range("A1").Value = "targ_right"
$ws.Range("B1").Value = "cue_valid"

# New data rows - all 4 combinations of targ_right/cue_valid.
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 0

$ws.Range("A3").Value = 0
$ws.Range("B3").Value = 1

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 0

$ws.Range("A5").Value = 1
$ws.Range("B5").Value = 1

# Drop the now-unused columns (session_id, trial_id, TRIAL_START,
# TRIAL_END, corr_resp, rt) that used to live in A:D and G:H.
$ws.Range("C1:H5").ClearContents()
